$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 ("Marking"): Right mark value and Wrong mark value
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 ("Total"): Right total and Wrong total, plus the "x / y" summary text
$ws.Range("B12").Value = 68
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "66 / 112"
